$d = $word.ActiveDocument

# The bibliography section ends with "... Thomson Pioneira (2008)." and was
# then followed (in the old site build) by a blank paragraph plus two
# boilerplate paragraphs scraped from the page chrome:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
# This rebuild of the site dropped that chrome. Locate those two text
# paragraphs plus the blank paragraph immediately before them, and delete
# the whole block, leaving the trailing blank/page-break paragraphs intact.

$verIndex = -1
$copyrightIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") { $verIndex = $i }
    if ($t -like "*Contact: luizeleno@usp.br*") { $copyrightIndex = $i }
}

if ($verIndex -gt 0 -and $copyrightIndex -ge $verIndex) {
    $blankIndex = $verIndex - 1
    $startPara = $d.Paragraphs.Item($blankIndex)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
